# "start new week 12/2"
# Roll the weekly Agenda tracker forward by one week:
#  - shift the Mon..Sun dates in A2:A8 forward by 7 days
#  - the daily hour-tracking cells (C:L) for the new (unstarted) week reset to 0
#  - the last day row (row 8) picks up the same look (borders/number format) as
#    the rest of the date column, so copy the format down from row 7
#  - the totals (row 9 SUM, row 11 diff, B13 SUM) recalc automatically
#  - the view is left focused on the date column that was just updated, at a
#    slightly smaller zoom

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agenda")

# --- shift the week's dates forward by 7 days -----------------------------
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value + 7
}

# --- new week -> no hours logged yet, reset the daily tracker cells -------
$ws.Range("C2:L8").Value = 0

# --- row 8 (Sun) should look like the rest of the date column again -------
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- restore the exact numeric/date values after the format paste ---------
$ws.Cells.Item(2, 1).Value = 43801
$ws.Cells.Item(3, 1).Value = 43802
$ws.Cells.Item(4, 1).Value = 43803
$ws.Cells.Item(5, 1).Value = 43804
$ws.Cells.Item(6, 1).Value = 43805
$ws.Cells.Item(7, 1).Value = 43806
$ws.Cells.Item(8, 1).Value = 43807

# --- view: zoom out a little and select the date range just edited --------
$ws.Range("A2:A8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120
